$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 283; this shifts existing rows 283-296 down to 284-297
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with the new data record
$ws.Cells.Item(283, 1).Value = 9
$ws.Cells.Item(283, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(283, 3).Value = "Metropolitana"
$ws.Cells.Item(283, 4).Value = 44753
$ws.Cells.Item(283, 5).Value = 13
$ws.Cells.Item(283, 6).Value = 100112043
$ws.Cells.Item(283, 7).Value = "Pepino ensalada"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 70
$ws.Cells.Item(283, 11).Value = 16000
$ws.Cells.Item(283, 12).Value = 18000
$ws.Cells.Item(283, 13).Value = 17000
$ws.Cells.Item(283, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(283, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(283, 16).Value = 283
$ws.Cells.Item(283, 17).Value = 60
$ws.Cells.Item(283, 18).Value = "Hortaliza"
